$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fgf9"
$ws.Cells.Item(2, 3).Value = "Fgfr3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.7364236666666667
$ws.Cells.Item(2, 8).Value = 2.209271
$ws.Cells.Item(2, 9).Value = 0.9917500467982164
$ws.Cells.Item(2, 10).Value = 0.9944848646626661
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.330840333333333
$ws.Cells.Item(2, 14).Value = 6.992521
$ws.Cells.Item(2, 15).Value = 0.6715345129768794
$ws.Cells.Item(2, 16).Value = 0.7003397275969581
$ws.Cells.Item(2, 17).Value = 1.716485984687889
$ws.Cells.Item(2, 18).Value = 15.448373862191
$ws.Cells.Item(2, 19).Value = 0.6659943846714376
$ws.Cells.Item(2, 20).Value = 0.6964772592171492

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fgf9"
$ws.Cells.Item(3, 3).Value = "Fgfr3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.7364236666666667
$ws.Cells.Item(3, 8).Value = 2.209271
$ws.Cells.Item(3, 9).Value = 0.9917500467982164
$ws.Cells.Item(3, 10).Value = 0.9944848646626661
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.657666
$ws.Cells.Item(3, 14).Value = 1.972998
$ws.Cells.Item(3, 15).Value = 0.1894790521235985
$ws.Cells.Item(3, 16).Value = 0.1976066831789769
$ws.Cells.Item(3, 17).Value = 0.484320807162
$ws.Cells.Item(3, 18).Value = 4.358887264458001
$ws.Cells.Item(3, 19).Value = 0.1879158588108605
$ws.Cells.Item(3, 20).Value = 0.1965168555776832

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fgf9"
$ws.Cells.Item(4, 3).Value = "Fgfr3"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.7364236666666667
$ws.Cells.Item(4, 8).Value = 2.209271
$ws.Cells.Item(4, 9).Value = 0.9917500467982164
$ws.Cells.Item(4, 10).Value = 0.9944848646626661
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.42828
$ws.Cells.Item(4, 14).Value = 0.85656
$ws.Cells.Item(4, 15).Value = 0.123391035029171
$ws.Cells.Item(4, 16).Value = 0.08578923067523865
$ws.Cells.Item(4, 17).Value = 0.31539552796
$ws.Cells.Item(4, 18).Value = 1.89237316776
$ws.Cells.Item(4, 19).Value = 0.1223730647646607
$ws.Cells.Item(4, 20).Value = 0.08531609145757894

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Fgf9"
$ws.Cells.Item(5, 3).Value = "Fgfr3"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.7364236666666667
$ws.Cells.Item(5, 8).Value = 2.209271
$ws.Cells.Item(5, 9).Value = 0.9917500467982164
$ws.Cells.Item(5, 10).Value = 0.9944848646626661
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.05413033333333334
$ws.Cells.Item(5, 14).Value = 0.162391
$ws.Cells.Item(5, 15).Value = 0.01559539987035126
$ws.Cells.Item(5, 16).Value = 0.01626435854882633
$ws.Cells.Item(5, 17).Value = 0.03986285855122223
$ws.Cells.Item(5, 18).Value = 0.3587657269610001
$ws.Cells.Item(5, 19).Value = 0.01546673855125776
$ws.Cells.Item(5, 20).Value = 0.01617465841025463

# Row 6
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Fgf9"
$ws.Cells.Item(6, 3).Value = "Fgfr3"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.5
$ws.Cells.Item(6, 7).Value = 0.006126
$ws.Cells.Item(6, 8).Value = 0.012252
$ws.Cells.Item(6, 9).Value = 0.008249953201783585
$ws.Cells.Item(6, 10).Value = 0.005515135337333892
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.330840333333333
$ws.Cells.Item(6, 14).Value = 6.992521
$ws.Cells.Item(6, 15).Value = 0.6715345129768794
$ws.Cells.Item(6, 16).Value = 0.7003397275969581
$ws.Cells.Item(6, 17).Value = 0.014278727882
$ws.Cells.Item(6, 18).Value = 0.085672367292
$ws.Cells.Item(6, 19).Value = 0.005540128305441787
$ws.Cells.Item(6, 20).Value = 0.003862468379808775

# Row 7
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Fgf9"
$ws.Cells.Item(7, 3).Value = "Fgfr3"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.5
$ws.Cells.Item(7, 7).Value = 0.006126
$ws.Cells.Item(7, 8).Value = 0.012252
$ws.Cells.Item(7, 9).Value = 0.008249953201783585
$ws.Cells.Item(7, 10).Value = 0.005515135337333892
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.657666
$ws.Cells.Item(7, 14).Value = 1.972998
$ws.Cells.Item(7, 15).Value = 0.1894790521235985
$ws.Cells.Item(7, 16).Value = 0.1976066831789769
$ws.Cells.Item(7, 17).Value = 0.004028861916
$ws.Cells.Item(7, 18).Value = 0.024173171496
$ws.Cells.Item(7, 19).Value = 0.001563193312738
$ws.Cells.Item(7, 20).Value = 0.001089827601293718

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Fgf9"
$ws.Cells.Item(8, 3).Value = "Fgfr3"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.5
$ws.Cells.Item(8, 7).Value = 0.006126
$ws.Cells.Item(8, 8).Value = 0.012252
$ws.Cells.Item(8, 9).Value = 0.008249953201783585
$ws.Cells.Item(8, 10).Value = 0.005515135337333892
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.42828
$ws.Cells.Item(8, 14).Value = 0.85656
$ws.Cells.Item(8, 15).Value = 0.123391035029171
$ws.Cells.Item(8, 16).Value = 0.08578923067523865
$ws.Cells.Item(8, 17).Value = 0.00262364328
$ws.Cells.Item(8, 18).Value = 0.01049457312
$ws.Cells.Item(8, 19).Value = 0.0010179702645103
$ws.Cells.Item(8, 20).Value = 0.0004731392176596974

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Fgf9"
$ws.Cells.Item(9, 3).Value = "Fgfr3"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.5
$ws.Cells.Item(9, 7).Value = 0.006126
$ws.Cells.Item(9, 8).Value = 0.012252
$ws.Cells.Item(9, 9).Value = 0.008249953201783585
$ws.Cells.Item(9, 10).Value = 0.005515135337333892
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.05413033333333334
$ws.Cells.Item(9, 14).Value = 0.162391
$ws.Cells.Item(9, 15).Value = 0.01559539987035126
$ws.Cells.Item(9, 16).Value = 0.01626435854882633
$ws.Cells.Item(9, 17).Value = 0.000331602422
$ws.Cells.Item(9, 18).Value = 0.001989614532
$ws.Cells.Item(9, 19).Value = 0.0001286613190934997
$ws.Cells.Item(9, 20).Value = 0.00008970013857170065

